$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 7, shifting existing rows 7-26 down to 8-27.
# Excel copies the formatting (e.g. the date number format) of the row above
# into the newly inserted row, matching style s="2" already used by column D.
$ws.Rows.Item(7).Insert()

# Populate the new weekly record in row 7 (new date 2022-01-19 / serial 44580).
$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "Vega Monumental Concepción"
$ws.Range("C7").Value = "Bíobío"
$ws.Range("D7").Value = 44580
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 100112030
$ws.Range("G7").Value = "Poroto granado"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 28000
$ws.Range("L7").Value = 30000
$ws.Range("M7").Value = 29000
$ws.Range("N7").Value = "$/saco 25 kilos"
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value = 1160
$ws.Range("Q7").Value = 25
$ws.Range("R7").Value = "Hortaliza"
